$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Gallery - A page of photos" list item entirely (including
#    its own paragraph mark), merging it away so the following "Sources -
#    Where to find out about the person" item moves up into its place.
# ---------------------------------------------------------------------------
$galleryRng = $d.Content
$galleryRng.Find.ClearFormatting()
$galleryFound = $galleryRng.Find.Execute("Gallery")
if ($galleryFound) {
  $galleryPara = $galleryRng.Paragraphs(1)
  $nextPara = $galleryPara.Next()
  $killRng = $d.Range($galleryPara.Range.Start, $nextPara.Range.Start)
  $killRng.Delete()
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark (Word's "last edit position" marker) to sit
#    right at the start of the paragraph that now begins with "Sources"
#    (this is simply where the previous edit left the cursor).
# ---------------------------------------------------------------------------
$sourcesRng = $d.Content
$sourcesRng.Find.ClearFormatting()
$sourcesFound = $sourcesRng.Find.Execute("Sources")
if ($sourcesFound) {
  $sourcesPara = $sourcesRng.Paragraphs(1)
  $pos = $sourcesPara.Range.Start
  try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
  } catch {
  }
  $d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos)) | Out-Null
}

# ---------------------------------------------------------------------------
# Shared run-properties fragment used by the reconstructed runs below.
# ---------------------------------------------------------------------------
$rpr = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# ---------------------------------------------------------------------------
# 3) Split ". (please use .docx " into three runs around "docx" and wrap
#    "docx" with spell-check proofErr tags, leaving the following
#    "or .pdf " / "format)" runs untouched.
# ---------------------------------------------------------------------------
$docxRng = $d.Content
$docxRng.Find.ClearFormatting()
$docxFound = $docxRng.Find.Execute(". (please use .docx ")
if ($docxFound) {
  $docxPara = $docxRng.Paragraphs(1)
  $tailRng = $d.Range($docxRng.Start, $docxPara.Range.End - 1)

  $inner = "<w:r w:rsidRPr=`"007C4EF5`">$rpr<w:t>. (please use .</w:t></w:r>" +
           "<w:proofErr w:type=`"spellStart`"/>" +
           "<w:r>$rpr<w:t>docx</w:t></w:r>" +
           "<w:proofErr w:type=`"spellEnd`"/>" +
           "<w:r>$rpr<w:t xml:space=`"preserve`"> </w:t></w:r>" +
           "<w:r>$rpr<w:t xml:space=`"preserve`">or .pdf </w:t></w:r>" +
           "<w:r w:rsidRPr=`"007C4EF5`">$rpr<w:t>format)</w:t></w:r>"

  $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
             '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
             '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
             '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
             $inner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

  $tailRng.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 4) Split "Items 1 and 2 above, but revised as needed." into three runs,
#    wrapping "above, but" with grammar-check proofErr tags.
# ---------------------------------------------------------------------------
$itemsRng = $d.Content
$itemsRng.Find.ClearFormatting()
$itemsFound = $itemsRng.Find.Execute("Items 1 and 2 above, but revised as needed.")
if ($itemsFound) {
  $itemsPara = $itemsRng.Paragraphs(1)
  $tailRng2 = $d.Range($itemsRng.Start, $itemsPara.Range.End - 1)

  $inner2 = "<w:r>$rpr<w:t xml:space=`"preserve`">Items 1 and 2 </w:t></w:r>" +
            "<w:proofErr w:type=`"gramStart`"/>" +
            "<w:r>$rpr<w:t>above, but</w:t></w:r>" +
            "<w:proofErr w:type=`"gramEnd`"/>" +
            "<w:r>$rpr<w:t xml:space=`"preserve`"> revised as needed.</w:t></w:r>"

  $xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
              '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
              '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
              '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
              $inner2 + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

  $tailRng2.InsertXML($xmlFrag2)
}

Write-Output "done"
